# Insert a new row for the "COMPLETED" / "Zdobyte:" / "Completed:" translation
# entry right before the existing "ANSWER50QUESTIONS" row (row 34), shifting
# all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 34..43 down to 35..44, leaving row 34 empty.
$ws.Rows("34").Insert()

# Populate the new row 34 with the new translation entry.
$ws.Range("A34").Value = "COMPLETED"
$ws.Range("B34").Value = "Zdobyte:"
$ws.Range("C34").Value = "Completed:"

# Update the active selection to match the new layout.
$ws.Range("C34").Select()
